# Auto-generated COM-interop script applying the meteocat daily-summary
# refresh described in the commit "Update automàtic: dades i banners
# [2026-02-19 20:50]". Each statement rewrites one cell's text to the
# freshly re-extracted value (extraction timestamp and/or measurement).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-19 20:48:37"
$ws.Range("I2").Value = "3.3 mm"
$ws.Range("E3").Value = "2026-02-19 20:48:40"
$ws.Range("I3").Value = "5.0 mm"
$ws.Range("E4").Value = "2026-02-19 20:48:43"
$ws.Range("J4").Value = "1009.9 hPa"
$ws.Range("E5").Value = "2026-02-19 20:48:46"
$ws.Range("I5").Value = "7.6 mm"
$ws.Range("E6").Value = "2026-02-19 20:48:49"
$ws.Range("J6").Value = "1010.0 hPa"
$ws.Range("O6").Value = "10.4 °C"
$ws.Range("E7").Value = "2026-02-19 20:48:51"
$ws.Range("J7").Value = "1011.0 hPa"
$ws.Range("E8").Value = "2026-02-19 20:48:54"
$ws.Range("J8").Value = "1010.7 hPa"
$ws.Range("O8").Value = "9.9 °C"
$ws.Range("E9").Value = "2026-02-19 20:48:57"
$ws.Range("E10").Value = "2026-02-19 20:48:59"
$ws.Range("N10").Value = "4.3 °C 20:23 TU"
$ws.Range("O10").Value = "10.4 °C"
$ws.Range("E11").Value = "2026-02-19 20:49:02"
$fmt = $ws.Range("H11").NumberFormat
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "64%"
$ws.Range("H11").NumberFormat = $fmt
$ws.Range("O11").Value = "5.7 °C"
$ws.Range("E12").Value = "2026-02-19 20:49:04"
$fmt = $ws.Range("H12").NumberFormat
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "78%"
$ws.Range("H12").NumberFormat = $fmt
$ws.Range("E13").Value = "2026-02-19 20:49:07"
$fmt = $ws.Range("H13").NumberFormat
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "62%"
$ws.Range("H13").NumberFormat = $fmt
$ws.Range("O13").Value = "4.4 °C"
$ws.Range("E14").Value = "2026-02-19 20:49:10"
$ws.Range("E15").Value = "2026-02-19 20:49:12"
$ws.Range("E16").Value = "2026-02-19 20:49:13"
$ws.Range("I16").Value = "9.7 mm"
$ws.Range("E17").Value = "2026-02-19 20:49:15"
$ws.Range("E18").Value = "2026-02-19 20:49:16"
$ws.Range("J18").Value = "1010.2 hPa"
$ws.Range("O18").Value = "11.7 °C"
$ws.Range("E19").Value = "2026-02-19 20:49:18"
$ws.Range("O19").Value = "5.3 °C"
$ws.Range("E20").Value = "2026-02-19 20:49:21"
$ws.Range("E21").Value = "2026-02-19 20:49:24"
$fmt = $ws.Range("H21").NumberFormat
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "61%"
$ws.Range("H21").NumberFormat = $fmt
$ws.Range("J21").Value = "1011.3 hPa"
$ws.Range("E22").Value = "2026-02-19 20:49:25"
$ws.Range("I22").Value = "0.3 mm"
$ws.Range("L22").Value = "99.7 km/h - 343º 20:23 TU"
$ws.Range("E23").Value = "2026-02-19 20:49:28"
$ws.Range("I23").Value = "10.1 mm"
$ws.Range("O23").Value = "-6.5 °C"
$ws.Range("E24").Value = "2026-02-19 20:49:30"
$ws.Range("J24").Value = "1014.9 hPa"
$ws.Range("E25").Value = "2026-02-19 20:49:33"
$ws.Range("I25").Value = "6.4 mm"
$ws.Range("E26").Value = "2026-02-19 20:49:36"
$ws.Range("J26").Value = "1009.9 hPa"
$ws.Range("L26").Value = "70.2 km/h - 307º 20:02 TU"
$ws.Range("E27").Value = "2026-02-19 20:49:39"
$ws.Range("O27").Value = "-3.7 °C"
$ws.Range("E28").Value = "2026-02-19 20:49:42"
$ws.Range("J28").Value = "1009.8 hPa"
$ws.Range("E29").Value = "2026-02-19 20:49:44"
$ws.Range("N29").Value = "5.4 °C 20:14 TU"
$ws.Range("O29").Value = "10.6 °C"
$ws.Range("E30").Value = "2026-02-19 20:49:47"
$ws.Range("J30").Value = "1010.0 hPa"
$ws.Range("O30").Value = "10.0 °C"
$ws.Range("E31").Value = "2026-02-19 20:49:49"
$ws.Range("J31").Value = "1009.5 hPa"
$ws.Range("E32").Value = "2026-02-19 20:49:52"
$ws.Range("E33").Value = "2026-02-19 20:49:54"
$ws.Range("J33").Value = "1010.8 hPa"
$ws.Range("E34").Value = "2026-02-19 20:49:57"
$ws.Range("E35").Value = "2026-02-19 20:50:00"
$ws.Range("J35").Value = "1016.3 hPa"
$ws.Range("E36").Value = "2026-02-19 20:50:03"
$ws.Range("J36").Value = "1010.3 hPa"
$ws.Range("E37").Value = "2026-02-19 20:50:05"
$ws.Range("J37").Value = "1011.3 hPa"
$ws.Range("E38").Value = "2026-02-19 20:50:08"
$ws.Range("E39").Value = "2026-02-19 20:50:11"
$ws.Range("I39").Value = "4.9 mm"
$ws.Range("E40").Value = "2026-02-19 20:50:14"
$fmt = $ws.Range("H40").NumberFormat
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "71%"
$ws.Range("H40").NumberFormat = $fmt
$ws.Range("J40").Value = "1012.5 hPa"
$ws.Range("E41").Value = "2026-02-19 20:50:16"
$ws.Range("J41").Value = "1012.9 hPa"
$ws.Range("O41").Value = "14.1 °C"
$ws.Range("E42").Value = "2026-02-19 20:50:19"
$ws.Range("O42").Value = "11.2 °C"
$ws.Range("E43").Value = "2026-02-19 20:50:22"
$ws.Range("E44").Value = "2026-02-19 20:50:24"
$ws.Range("E45").Value = "2026-02-19 20:50:27"
$ws.Range("J45").Value = "1015.5 hPa"
$ws.Range("O45").Value = "2.6 °C"
$ws.Range("E46").Value = "2026-02-19 20:50:29"
$ws.Range("J46").Value = "1015.7 hPa"
